$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 4

$ws.Cells.Item($row, 1).Value = "2025-11-07 05:41:19 UTC"
$ws.Cells.Item($row, 2).Value = "2025-11-07 11:11:19 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-11-2025.pdf"
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = ""

$ws.Range("A3:H3").Copy()
$ws.Range("A4:H4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
